$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H105").Value = 59223.332
$ws.Range("J105").Value = 59223.332
$ws.Range("L105").Value = 59223.332
$ws.Range("N105").Value = -66211.33199999999
$ws.Range("H109").Value = 60684
$ws.Range("J109").Value = 60684
$ws.Range("L109").Value = 60684
$ws.Range("N109").Value = -63458
$ws.Range("H111").Value = 2530.5
$ws.Range("I111").Value = 2029
$ws.Range("J111").Value = 3032
$ws.Range("K111").Value = 6087
$ws.Range("L111").Value = 9096
$ws.Range("M111").Value = -3020
$ws.Range("N111").Value = -15230
$ws.Range("H113").Value = 26938.125
$ws.Range("J113").Value = 3000
$ws.Range("L113").Value = 3000
$ws.Range("N113").Value = -9508
$ws.Range("H114").Value = 60722
$ws.Range("J114").Value = 60722
$ws.Range("L114").Value = 60722
$ws.Range("N114").Value = -69400
$ws.Range("H131").Value = 1448.15
$ws.Range("I131").Value = 652.0909
$ws.Range("J131").Value = 2421.111
$ws.Range("K131").Value = 1956.2727
$ws.Range("L131").Value = 7263.333
$ws.Range("M131").Value = 3083.7273
$ws.Range("N131").Value = -17343.333
$ws.Range("H132").Value = 1065.3928
$ws.Range("I132").Value = 955.4091
$ws.Range("K132").Value = 2866.2273
$ws.Range("M132").Value = -336.2273
$ws.Range("H137").Value = 1396.0667
$ws.Range("I137").Value = 995.46155
$ws.Range("K137").Value = 2986.38465
$ws.Range("M137").Value = -436.38465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1389926.2
$ws.Range("I2").Value = 2777853
$ws.Range("K2").Value = 2777853
$ws.Range("M2").Value = -2777740
$ws.Range("H32").Value = 3620.125
$ws.Range("I32").Value = 2970
$ws.Range("K32").Value = 2970
$ws.Range("M32").Value = -2683
$ws.Range("H45").Value = 1665.4546
$ws.Range("I45").Value = 1056.5
$ws.Range("K45").Value = 1056.5
$ws.Range("M45").Value = -679.5
$ws.Range("H74").Value = 1766.6538
$ws.Range("I74").Value = 458.17648
$ws.Range("K74").Value = 458.17648
$ws.Range("M74").Value = 415.82352
$ws.Range("H77").Value = 1766.6538
$ws.Range("I77").Value = 458.17648
$ws.Range("K77").Value = 2290.8824
$ws.Range("M77").Value = 2077.1176
$ws.Range("H116").Value = 1389926.2
$ws.Range("I116").Value = 2777853
$ws.Range("K116").Value = 2777853
$ws.Range("M116").Value = -2775559
$ws.Range("H123").Value = 74250
$ws.Range("J123").Value = 74250
$ws.Range("L123").Value = 74250
$ws.Range("N123").Value = -84050
$ws.Range("H132").Value = 1262.875
$ws.Range("I132").Value = 1135.5676
$ws.Range("K132").Value = 3406.7028
$ws.Range("M132").Value = -876.7028

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1389926.2
$ws.Range("I3").Value = 2777853
$ws.Range("K3").Value = 2777853
$ws.Range("M3").Value = -2777739
$ws.Range("H134").Value = 2210.3333
$ws.Range("I134").Value = 1966.6923
$ws.Range("J134").Value = 2606.25
$ws.Range("K134").Value = 5900.0769
$ws.Range("L134").Value = 7818.75
$ws.Range("M134").Value = -3365.0769
$ws.Range("N134").Value = -12888.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 716.8570999999999
$ws.Range("I16").Value = 758.8
$ws.Range("K16").Value = 758.8
$ws.Range("M16").Value = -471.8
$ws.Range("H31").Value = 1980.6296
$ws.Range("I31").Value = 1666.3
$ws.Range("J31").Value = 2165.5293
$ws.Range("K31").Value = 1666.3
$ws.Range("L31").Value = 2165.5293
$ws.Range("M31").Value = -1371.3
$ws.Range("N31").Value = -2755.5293
$ws.Range("H34").Value = 1980.6296
$ws.Range("I34").Value = 1666.3
$ws.Range("J34").Value = 2165.5293
$ws.Range("K34").Value = 1666.3
$ws.Range("L34").Value = 2165.5293
$ws.Range("M34").Value = -1464.3
$ws.Range("N34").Value = -2569.5293
$ws.Range("H105").Value = 1142.7142
$ws.Range("I105").Value = 1148
$ws.Range("K105").Value = 1148
$ws.Range("M105").Value = 599
$ws.Range("H113").Value = 716.8570999999999
$ws.Range("I113").Value = 758.8
$ws.Range("K113").Value = 758.8
$ws.Range("M113").Value = 1411.2
$ws.Range("H141").Value = 65494.5
$ws.Range("J141").Value = 63326
$ws.Range("L141").Value = 63326
$ws.Range("N141").Value = -73686

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1456.4286
$ws.Range("I81").Value = 822.5
$ws.Range("K81").Value = 2467.5
$ws.Range("M81").Value = -1344.5
$ws.Range("H84").Value = 1456.4286
$ws.Range("I84").Value = 822.5
$ws.Range("K84").Value = 7402.5
$ws.Range("M84").Value = -1786.5
$ws.Range("H131").Value = 17829.586
$ws.Range("J131").Value = 19680.621
$ws.Range("L131").Value = 59041.863
$ws.Range("N131").Value = -69121.863

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2607.1765
$ws.Range("I7").Value = 2488.875
$ws.Range("J7").Value = 4500
$ws.Range("K7").Value = 2488.875
$ws.Range("L7").Value = 4500
$ws.Range("M7").Value = -2376.875
$ws.Range("N7").Value = -4724
$ws.Range("H46").Value = 1879.9231
$ws.Range("I46").Value = 1111.4286
$ws.Range("J46").Value = 2776.5
$ws.Range("K46").Value = 1111.4286
$ws.Range("L46").Value = 2776.5
$ws.Range("M46").Value = -923.4286
$ws.Range("N46").Value = -3152.5
$ws.Range("H126").Value = 2607.1765
$ws.Range("I126").Value = 2488.875
$ws.Range("J126").Value = 4500
$ws.Range("K126").Value = 7466.625
$ws.Range("L126").Value = 13500
$ws.Range("M126").Value = -4996.625
$ws.Range("N126").Value = -18440

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7754.1816
$ws.Range("J96").Value = 11971
$ws.Range("L96").Value = 11971
$ws.Range("N96").Value = -14717
$ws.Range("H126").Value = 8524.799999999999
$ws.Range("I126").Value = 11035.363
$ws.Range("J126").Value = 1620.75
$ws.Range("K126").Value = 33106.089
$ws.Range("L126").Value = 4862.25
$ws.Range("M126").Value = -30636.089
$ws.Range("N126").Value = -9802.25
$ws.Range("H132").Value = 1950.1562
$ws.Range("I132").Value = 1388
$ws.Range("J132").Value = 3186.9
$ws.Range("K132").Value = 4164
$ws.Range("L132").Value = 9560.700000000001
$ws.Range("M132").Value = -1634
$ws.Range("N132").Value = -14620.7
